$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 9 with the same pattern as the existing rows (Giovanni / Interno / 74 / 2019-01-30)
$ws.Range("A9").Value = "Giovanni"
$ws.Range("B9").Value = "Interno"
$ws.Range("C9").Value = 74

# Copy the date cell above (D8) so the new D9 cell picks up the exact same style (numFmtId 14)
$ws.Range("D8").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("D9").Value = 43495

# Move the active selection to C10, matching the post-edit state
$ws.Range("C10").Select()
